$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark from its old location
#    (right after the word "had"). Deleting a bookmark does not
#    remove/shift any text, it only drops the bookmark markers.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Grab a tiny "template" range elsewhere in the doc that already
#    carries the exact run formatting we need to reproduce
#    (color 000000 / sz 18 / szCs 18, nothing else) so that newly
#    inserted runs get a matching <w:rPr> including szCs, which is
#    not reachable through Font.Size alone.
# ------------------------------------------------------------------
$tmplFind = $d.Content
$tmplFind.Find.Execute("All files", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tmpl = $tmplFind.Duplicate
$tmpl.Start = $tmplFind.Start
$tmpl.End = $tmplFind.Start + 1
$tmplFormatted = $tmpl.FormattedText

# ------------------------------------------------------------------
# 3) Locate the sentence that needs to be rewritten / split into
#    several runs.
# ------------------------------------------------------------------
$old = "Each lock has a time to live (default 10 minutes) If I client ties to request a file that is locked, it will get a message saying that the file is locked and the TTL. Once the TTL expires, the lock will be released."
$target = $d.Content
$target.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startPos = $target.Start

# First chunk keeps the original run, just with shortened text.
$chunk1 = "Each lock has a time to live (default 10 minutes)"
$target.Text = $chunk1
$pos = $startPos + $chunk1.Length

function Insert-Chunk($atPos, $text) {
    $r = $d.Range($atPos, $atPos)
    $r.FormattedText = $tmplFormatted
    $fresh = $d.Range($atPos, $atPos + 1)
    $fresh.Text = $text
    return ($atPos + $text.Length)
}

$pos = Insert-Chunk $pos "."
$pos = Insert-Chunk $pos " If "
$pos = Insert-Chunk $pos "the "
$pos = Insert-Chunk $pos "client "
$pos = Insert-Chunk $pos "tries "

# ------------------------------------------------------------------
# 4) Re-insert the "_GoBack" bookmark (zero length) right before the
#    final chunk.
# ------------------------------------------------------------------
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$pos = Insert-Chunk $pos "to request a file that is locked, it will get a message saying that the file is locked and the TTL. Once the TTL expires, the lock will be released."
